$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: LUNETTES DE SOLEIL RANDONNEE
$ws.Range("A16").Value = "LUNETTES DE SOLEIL RANDONNÉE - MH140 - ADULTE - POLARISANTES CATÉGORIE 3"
$ws.Range("B16").Value = "C:\SymfonyProjects\e-commerce\public\pictures\lunettes-de-soleil-randonnee-mh140-adulte-polarisantes-categorie-3.avif"
$ws.Range("C16").Value = "Nos ingénieurs optiques ont développé ces lunettes de soleil pour la randonnée. Idéales pour un usage occasionnel en montagne grâce à leur légèreté. Les verres anti-UV bloquent 100% des rayons nocifs et la catégorie 3 vous protège de l’éblouissement. La technologie des verres polarisants vous permet de mieux distinguer les reliefs et contrastes."

# Row 17: CHAUSSETTES DE RUNNING
$ws.Range("A17").Value = "CHAUSSETTES DE RUNNING RUN100 NOIRE X3"
$ws.Range("B17").Value = "C:\SymfonyProjects\e-commerce\public\pictures\chaussettes-de-running-run100-noire-x3.avif"
$ws.Range("C17").Value = "Nos équipes de conception ont développé ces chaussettes de running pour que votre pied soit protégé des risques d'ampoules lors de vos sorties de course à pied. Des chaussettes de running à petit prix ? Grâce à leur fil en coton, leur tricotage fin et aéré, elles sont idéales pour la course à pied à petit budget."

# Row 18: MONTRE CARDIO GPS GARMIN
$ws.Range("A18").Value = "MONTRE CARDIO GPS GARMIN FORERUNNER 245 GRISE"
$ws.Range("B18").Value = "C:\SymfonyProjects\e-commerce\public\pictures\montre-cardio-gps-garmin-forerunner-245-grise.avif"
$ws.Range("C18").Value = "Conçue pour les runners, coureurs de 10km, de semi et de marathon. Elle suit vos statistiques, traite les données et acquiert toutes les informations sur vos performances, votre technique de course, l'historique de vos entraînements et même vos objectifs."

# Row 19: SAC A DOS DE RANDONNEE (path registered before name in shared strings)
$ws.Range("B19").Value = "C:\SymfonyProjects\e-commerce\public\pictures\sac-a-dos-de-randonnee-30l-nh-arpenaz-500.avif"
$ws.Range("A19").Value = "SAC À DOS DE RANDONNÉE 30L - NH ARPENAZ 500"
$ws.Range("C19").Value = "Nos concepteurs randonneurs ont conçu ce sac à dos NH Arpenaz 500 30 litres pour accompagner vos randonnées à la journée en plaine, forêt ou sur le littoral. Notre motivation ? Vous proposer un sac à dos confortable et très accessoirisé pour profiter de vos randonnées ! Retrouvez une poche pour conserver au frais votre pique-nique et une poche téléphone."

# Widen column B (and nudge column C) to fit the new, longer content
# (column A is left as-is: its target width is within COM rounding noise of the original)
$ws.Columns.Item(2).ColumnWidth = 162.3333
$ws.Columns.Item(3).ColumnWidth = 157.4615

# Update selection to match target state
$ws.Range("B21").Select()
